$d = $word.ActiveDocument

$replacements = @(
    @{old="2025-05-19 Monday"; new="2025-05-20 Tuesday"},
    @{old="29÷6=4, 5"; new="55÷6=9, 1"},
    @{old="30÷8=3, 6"; new="74÷9=8, 2"},
    @{old="82÷5=16, 2"; new="39÷9=4, 3"},
    @{old="63÷9=7, 0"; new="22÷9=2, 4"},
    @{old="92÷6=15, 2"; new="22÷9=2, 4"},
    @{old="25÷9=2, 7"; new="32÷6=5, 2"},
    @{old="56÷6=9, 2"; new="33÷4=8, 1"},
    @{old="29÷7=4, 1"; new="70÷2=35, 0"},
    @{old="92÷8=11, 4"; new="38÷2=19, 0"},
    @{old="57÷9=6, 3"; new="23÷2=11, 1"},
    @{old="31÷4=7, 3"; new="37÷8=4, 5"},
    @{old="24÷6=4, 0"; new="85÷4=21, 1"},
    @{old="59÷6=9, 5"; new="46÷5=9, 1"},
    @{old="21÷7=3, 0"; new="56÷3=18, 2"},
    @{old="65÷6=10, 5"; new="53÷5=10, 3"},
    @{old="31÷7=4, 3"; new="51÷9=5, 6"},
    @{old="22÷6=3, 4"; new="47÷6=7, 5"},
    @{old="79÷2=39, 1"; new="73÷4=18, 1"},
    @{old="62÷6=10, 2"; new="11÷3=3, 2"},
    @{old="81÷4=20, 1"; new="95÷6=15, 5"},
    @{old="20÷4=5, 0"; new="13÷4=3, 1"},
    @{old="85÷7=12, 1"; new="65÷4=16, 1"},
    @{old="27÷7=3, 6"; new="72÷2=36, 0"},
    @{old="56÷8=7, 0"; new="34÷5=6, 4"},
    @{old="10÷3=3, 1"; new="23÷8=2, 7"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $r.new, 2)
}
